# Build.xlsx update — add "icon" (按钮图片) and "texture" (贴图) columns
# to the build table, and fix the Build_Bakery level_need value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (row 1: display name, row 2: field key, row 3: field type) ---
$ws.Range("E1").Value = "按钮图片"
$ws.Range("F1").Value = "贴图"

$ws.Range("E2").Value = "icon"
$ws.Range("F2").Value = "texture"

$ws.Range("E3").Value = "string"
$ws.Range("F3").Value = "string"

# --- Row 4: Build_Bakery — level_need corrected from 2 to 1, icon/texture added ---
$ws.Range("D4").Value = "1"
$ws.Range("E4").Value = "Build_Bakery"
$ws.Range("F4").Value = "Buid_01"

# --- Row 5: Build_Jam — icon added (no texture) ---
$ws.Range("E5").Value = "Build_Jam"

# --- Row 6: Build_Dairy — icon added (no texture) ---
$ws.Range("E6").Value = "Build_Dairy"

# --- Column widths for the two new columns ---
$ws.Columns("E").ColumnWidth = 14.8
$ws.Columns("F").ColumnWidth = 16.55

# --- Selection moves to D5 ---
$ws.Range("D5").Select()
